$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("drop")

# Column B ("n_obs" values) switches from numbers to text, row by row
$ws.Range("B2").Value = "'1"
$ws.Range("B2").Style = "Normal"
$ws.Range("B3").Value = "'0"
$ws.Range("B3").Style = "Normal"

# Column C mirrors column B
$ws.Range("C2").Value = "'22"
$ws.Range("C2").Style = "Normal"
$ws.Range("C3").Value = "'0"
$ws.Range("C3").Style = "Normal"

# Column D becomes empty text
$ws.Range("D2").Value = "'"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

# Column E becomes empty text
$ws.Range("E2").Value = "'"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'"
$ws.Range("E3").Style = "Normal"

# New column F: date_last_changed header + values
$ws.Range("F1").Value = "date_last_changed"
$ws.Range("F2").Value = "'11 Dec 2022"
$ws.Range("F2").Style = "Normal"
$ws.Range("F3").Value = "'11 Dec 2022"
$ws.Range("F3").Style = "Normal"
